$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 6: a second "Json Url request" row, this time linking to the
#     queryAttachments endpoint, with the URL rendered as a live hyperlink ---
$ws.Range("A6").Value = "Json Url request2"
$ws.Range("B6").Value = 'http://services5.arcgis.com/PnnKqtqi3qfxnaPc/arcgis/rest/services/LabOfflineGdbT3/FeatureServer/0/queryAttachments?objectIds=1169,1170,1171,1172,1173,1174,1175,1176,1177,1178,1179,1180,1181,1182,1183,1184&f=pjson&token='
[void]$ws.Hyperlinks.Add($ws.Range("B6"), 'http://services5.arcgis.com/PnnKqtqi3qfxnaPc/arcgis/rest/services/LabOfflineGdbT3/FeatureServer/0/queryAttachments?objectIds=1169,1170,1171,1172,1173,1174,1175,1176,1177,1178,1179,1180,1181,1182,1183,1184&f=pjson&token=')

# --- Row 8: the JSON result returned by that attachments query, wrapped
#     and given the same "max height" treatment as the row 4 JSON result ---
$ws.Range("A8").Value = "Json Result"
$ws.Range("B8").Value = '{
  "fields" : [
    {
      "name" : "OBJECTID", 
      "type" : "esriFieldTypeOID", 
      "alias" : "OBJECTID", 
      "sqlType" : "sqlTypeOther", 
      "domain" : null, 
      "defaultValue" : null
    }, 
    {
      "name" : "GlobalID", 
      "type" : "esriFieldTypeGlobalID", 
      "alias" : "GlobalID", 
      "sqlType" : "sqlTypeOther", 
      "length" : 38, 
      "domain" : null
    }
  ], 
  "attachmentGroups" : [
    {
      "parentObjectId" : 1169, 
      "parentGlobalId" : "f5f0fe49-0d01-4254-bf42-e1cb34d1187a", 
      "attachmentInfos" : [
        {
          "id" : 96, 
          "globalId" : "be7495cc-f392-46b1-940f-7121180627fc", 
          "name" : "Photo1.jpg", 
          "contentType" : "image/jpeg", 
          "size" : 206764, 
          "keywords" : ""
        }
      ]
    }, 
    {
      "parentObjectId" : 1170, 
      "parentGlobalId" : "4ed9c254-b8ec-420a-8a27-321d0af48866", 
      "attachmentInfos" : [
        {
          "id" : 97, 
          "globalId" : "434d0cae-a4c1-4999-87b4-1ee64cf606e3", 
          "name" : "Photo1.jpg", 
          "contentType" : "image/jpeg", 
          "size" : 609029, 
          "keywords" : ""
        }
      ]
    }, 
    {
      "parentObjectId" : 1171, 
      "parentGlobalId" : "4d214cfc-2626-499d-a747-2e8f2d216606", 
      "attachmentInfos" : [
        {
          "id" : 98, 
          "globalId" : "c4615073-c7df-46dd-9fef-4fcd05b644ae", 
          "name" : "Photo1.jpg", 
          "contentType" : "image/jpeg", 
          "size" : 648260, 
          "keywords" : ""
        }
      ]
    }, 
    {
      "parentObjectId" : 1172, 
      "parentGlobalId" : "a41fcf88-7c58-41ab-95f3-62340fa827a5", 
      "attachmentInfos" : [
        {
          "id" : 99, 
          "globalId" : "e01a005f-2ad4-4f58-9bef-6fa7e985a25c", 
          "name" : "Photo1.jpg", 
          "contentType" : "image/jpeg", 
          "size" : 588344, 
          "keywords" : ""
        }
      ]
    }, 
    {
      "parentObjectId" : 1173, 
      "parentGlobalId" : "ee811c17-4b96-4c12-afe5-111a59b535d8", 
      "attachmentInfos" : [
        {
          "id" : 100, 
          "globalId" : "b35a9e1c-f424-4c6c-9a84-332168e34036", 
          "name" : "Photo1.jpg", 
          "contentType" : "image/jpeg", 
          "size" : 650790, 
          "keywords" : ""
        }
      ]
    }, 
    {
      "parentObjectId" : 1176, 
      "parentGlobalId" : "8a45d50a-940f-431a-adc8-1e4f33944416", 
      "attachmentInfos" : [
        {
          "id" : 101, 
          "globalId" : "9da76763-bfc0-4ffa-a6c6-23e04a0f1a5b", 
          "name" : "Photo1.jpg", 
          "contentType" : "image/jpeg", 
          "size" : 603547, 
          "keywords" : ""
        }
      ]
    }, 
    {
      "parentObjectId" : 1177, 
      "parentGlobalId" : "29f39782-3bda-48e4-b0c5-f66732fc06e3", 
      "attachmentInfos" : [
        {
          "id" : 102, 
          "globalId" : "7239deaa-9182-41d4-be38-47843b97388a", 
          "name" : "Photo1.jpg", 
          "contentType" : "image/jpeg", 
          "size" : 647561, 
          "keywords" : ""
        }
      ]
    }, 
    {
      "parentObjectId" : 1178, 
      "parentGlobalId" : "37de1bde-67e5-449a-9e0e-85efc58d94eb", 
      "attachmentInfos" : [
        {
          "id" : 103, 
          "globalId" : "6ecdb19b-f52c-4273-9b84-81c921de9883", 
          "name" : "Photo1.jpg", 
          "contentType" : "image/jpeg", 
          "size" : 688391, 
          "keywords" : ""
        }
      ]
    }, 
    {
      "parentObjectId" : 1180, 
      "parentGlobalId" : "e5f579cf-4e0d-4b1a-b83e-a2ae867cb575", 
      "attachmentInfos" : [
        {
          "id" : 104, 
          "globalId" : "8e80c451-ffa0-4dfe-8021-afdfbed6f09b", 
          "name" : "Photo1.jpg", 
          "contentType" : "image/jpeg", 
          "size" : 526378, 
          "keywords" : ""
        }
      ]
    }, 
    {
      "parentObjectId" : 1181, 
      "parentGlobalId" : "5b8293ae-a562-466d-ba31-98994ebcdede", 
      "attachmentInfos" : [
        {
          "id" : 105, 
          "globalId" : "0b223829-d888-404e-9811-b60cf5709905", 
          "name" : "Photo1.jpg", 
          "contentType" : "image/jpeg", 
          "size" : 483106, 
          "keywords" : ""
        }
      ]
    }, 
    {
      "parentObjectId" : 1182, 
      "parentGlobalId" : "98225d2e-bb6c-4772-accc-0fc8be84206f", 
      "attachmentInfos" : [
        {
          "id" : 106, 
          "globalId" : "76f6fcdf-9c50-469a-be8e-253ed1f16b67", 
          "name" : "Photo1.jpg", 
          "contentType" : "image/jpeg", 
          "size" : 577698, 
          "keywords" : ""
        }
      ]
    }, 
    {
      "parentObjectId" : 1183, 
      "parentGlobalId" : "0bd696e2-68e8-4ccc-97b4-b963d5beb556", 
      "attachmentInfos" : [
        {
          "id" : 107, 
          "globalId" : "a4038103-4254-4767-a94d-0929bf3ca549", 
          "name" : "Photo1.jpg", 
          "contentType" : "image/jpeg", 
          "size" : 662714, 
          "keywords" : ""
        }
      ]
    }, 
    {
      "parentObjectId" : 1184, 
      "parentGlobalId" : "43b8bc07-d87d-4adf-8823-d6ee02fb11ec", 
      "attachmentInfos" : [
        {
          "id" : 108, 
          "globalId" : "a5f847c5-c8ed-446d-b12a-a056b835c690", 
          "name" : "Photo1.jpg", 
          "contentType" : "image/jpeg", 
          "size" : 670247, 
          "keywords" : ""
        }
      ]
    }
  ]
}'
$ws.Range("B8").WrapText = $true
$ws.Rows(8).RowHeight = 409.5

# Move the active selection down to the newly-added result cell, and (best
# effort) scroll the window so row 13 is pinned at the top, matching the
# saved view state of the authored workbook.
[void]$ws.Range("B8").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 13
$win.ScrollColumn = 1
